$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Force the Price column to text formatting before assignment so that
# Excel does not reinterpret decimal-looking strings (e.g. "0.999") as
# numbers, then restore the default "Normal" style so no visual/style
# change is introduced - only the literal text content changes.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '68.445.84'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +3.87%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.371.75'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +1.48%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '594.56'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +6.57%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '186.39'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.599'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +3.93%  '
$ws.Range('E9').Value = '  +4.79%  '
$ws.Range('E10').Value = '  +1.87%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '47.43'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +3.64%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0000282'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +7.41%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.915.12'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.37%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '639.72'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +12.15%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '8.59'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.58%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '68.610.70'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +3.93%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.387.06'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.58%  '
$ws.Range('E18').Value = '  +1.87%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '18.01'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.89%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.11'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +2.41%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.912'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.22%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '18.01'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('E23').Value = '  +1.86%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '99.45'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.09%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '4.10'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +3.73%  '
$ws.Range('E26').Value = '  +6.07%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.80'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +4.71%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '32.86'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +7.53%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.71'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +3.27%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.86'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.79%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '613.77'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +8.03%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.72'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.97%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.995.32'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +7.22%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '11.11'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.57%  '
$ws.Range('E35').Value = '  +2.89%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.28%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '56.28'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.34%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.80'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +8.21%  '
$ws.Range('E39').Value = '  +7.13%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.131'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +3.72%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '33.72'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.13%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0₃0709'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +3.04%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.43'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.32%  '
$ws.Range('E44').Value = '  +3.10%  '
$ws.Range('E45').Value = '  +4.03%  '
$ws.Range('E46').Value = '  +3.25%  '
$ws.Range('E47').Value = '  +3.73%  '
$ws.Range('E48').Value = '  +0.31%  '
$ws.Range('E49').Value = '  +11.61%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '131.61'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +5.19%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.81'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +7.22%  '
